$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (and E, for safety) to Text format while writing,
# so numeric-looking strings ("590.10", "0.0000265", ...) are not
# auto-converted to numbers / scientific notation by Excel.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "67.990.95"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "3.662.09"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "590.10"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "187.78"
$ws.Range("E6").Value = "  +3.17%  "

$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  -1.38%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "0.691"
$ws.Range("E9").Value = "  -3.26%  "

$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -7.34%  "

$ws.Range("D11").Value = "55.54"
$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("D12").Value = "0.0000265"
$ws.Range("E12").Value = "  -8.48%  "

$ws.Range("D13").Value = "10.13"
$ws.Range("E13").Value = "  -1.86%  "

$ws.Range("D14").Value = "4.242.37"
$ws.Range("E14").Value = "  -0.71%  "

$ws.Range("D15").Value = "3.656.83"
$ws.Range("E15").Value = "  -0.87%  "

$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "18.69"
$ws.Range("E17").Value = "  -2.97%  "

$ws.Range("D18").Value = "67.750.44"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("E19").Value = "  -2.42%  "

$ws.Range("D20").Value = "12.47"
$ws.Range("E20").Value = "  -2.26%  "

$ws.Range("D21").Value = "400.58"
$ws.Range("E21").Value = "  -1.88%  "

$ws.Range("E22").Value = "  -3.85%  "

$ws.Range("D23").Value = "87.09"
$ws.Range("E23").Value = "  -1.68%  "

$ws.Range("D24").Value = "2.91"
$ws.Range("E24").Value = "  -3.11%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "12.44"
$ws.Range("E25").Value = "  -2.47%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "10.71"
$ws.Range("E26").Value = "  -1.70%  "

$ws.Range("E27").Value = "  +0.74%  "

$ws.Range("D28").Value = "3.66"
$ws.Range("E28").Value = "  -4.94%  "

$ws.Range("D29").Value = "9.21"
$ws.Range("E29").Value = "  -1.97%  "

$ws.Range("D30").Value = "31.75"
$ws.Range("E30").Value = "  -2.86%  "

$ws.Range("D31").Value = "7.04"
$ws.Range("E31").Value = "  -2.85%  "

$ws.Range("D32").Value = "67.64"
$ws.Range("E32").Value = "  +5.37%  "

$ws.Range("D33").Value = "12.13"
$ws.Range("E33").Value = "  -2.32%  "

$ws.Range("D34").Value = "605.89"
$ws.Range("E34").Value = "  +1.10%  "

$ws.Range("D35").Value = "43.49"
$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").Value = "0.115"
$ws.Range("E36").Value = "  -1.79%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Value = "0.386"
$ws.Range("E39").Value = "  -2.92%  "

$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0748"
$ws.Range("E40").Value = "  -15.24%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.135"
$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").Value = "2.85"
$ws.Range("E42").Value = "  -4.86%  "

$ws.Range("D43").Value = "0.0419"
$ws.Range("E43").Value = "  -3.73%  "

$ws.Range("D44").Value = "2.48"
$ws.Range("E44").Value = "  -10.63%  "

$ws.Range("E45").Value = "  +0.65%  "

$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("D47").Value = "2.739.62"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("D48").Value = "8.74"
$ws.Range("E48").Value = "  -4.70%  "

$ws.Range("D49").Value = "141.93"
$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("D50").Value = "2.57"
$ws.Range("E50").Value = "  -5.41%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "2.66"
$ws.Range("E51").Value = "  -3.56%  "

# Restore the default (unstyled) cell style now that the text values
# are committed, so the cells stay visually/structurally unchanged.
$priceRange.Style = "Normal"
